# List Tugas Minggu Ke-10 — mark additional completed-task checkmarks.
# Cells that were blank now get the Wingdings "ü" (renders as a checkmark)
# used throughout the sheet to flag a finished task, matching the same
# font/fill formatting already used by the sibling checkmark cells in
# each column (C/E columns have no fill, the D column has a gray fill).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$checkMark = [char]252  # "ü" in Wingdings renders as a check mark glyph

function Set-Check($cellRef) {
    $rng = $ws.Range($cellRef)
    $rng.Value = $checkMark
    $rng.Font.Name = "Wingdings"
    $rng.HorizontalAlignment = -4108  # xlCenter
    $rng.VerticalAlignment = -4108    # xlCenter
}

# Row 9 — Dian Fitriani Aulia
Set-Check "C9"
Set-Check "D9"
Set-Check "E9"

# Row 10 — Dinda Fatchus Sabila Fa'izatur Rahmah
Set-Check "D10"

# Row 11 — Dinda Nur Apriani
Set-Check "D11"
Set-Check "E11"

# Row 14 — Keatryn Kezia P. Sihombing
Set-Check "D14"

# Row 24 — Nasywa Mawaddah
Set-Check "D24"
Set-Check "E24"

# Row 33 — Shabrina Cahyani
Set-Check "D33"

# Row 34 — Silvana Putri Ariani
Set-Check "D34"
Set-Check "E34"
